# Apply updated 'want to go' counts (column F) across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 37874
$ws.Cells.Item(4, 6).Value = 640
$ws.Cells.Item(6, 6).Value = 488
$ws.Cells.Item(7, 6).Value = 373
$ws.Cells.Item(8, 6).Value = 470
$ws.Cells.Item(9, 6).Value = 863
$ws.Cells.Item(10, 6).Value = 104
$ws.Cells.Item(11, 6).Value = 740
$ws.Cells.Item(12, 6).Value = 575
$ws.Cells.Item(13, 6).Value = 75
$ws.Cells.Item(16, 6).Value = 677
$ws.Cells.Item(17, 6).Value = 186
$ws.Cells.Item(18, 6).Value = 483
$ws.Cells.Item(19, 6).Value = 447
$ws.Cells.Item(20, 6).Value = 1184
$ws.Cells.Item(22, 6).Value = 862
$ws.Cells.Item(23, 6).Value = 2578
$ws.Cells.Item(24, 6).Value = 1063
$ws.Cells.Item(25, 6).Value = 576
$ws.Cells.Item(26, 6).Value = 113
$ws.Cells.Item(28, 6).Value = 46
$ws.Cells.Item(29, 6).Value = 817
$ws.Cells.Item(31, 6).Value = 1174

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 431
$ws.Cells.Item(4, 6).Value = 336
$ws.Cells.Item(10, 6).Value = 15

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 657

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 657
$ws.Cells.Item(3, 6).Value = 37874
$ws.Cells.Item(5, 6).Value = 640
$ws.Cells.Item(7, 6).Value = 488
$ws.Cells.Item(9, 6).Value = 373
$ws.Cells.Item(10, 6).Value = 470
$ws.Cells.Item(11, 6).Value = 431
$ws.Cells.Item(12, 6).Value = 336
$ws.Cells.Item(15, 6).Value = 863
$ws.Cells.Item(16, 6).Value = 104
$ws.Cells.Item(17, 6).Value = 740
$ws.Cells.Item(18, 6).Value = 575
$ws.Cells.Item(19, 6).Value = 75
$ws.Cells.Item(24, 6).Value = 15
$ws.Cells.Item(27, 6).Value = 677
$ws.Cells.Item(28, 6).Value = 186
$ws.Cells.Item(29, 6).Value = 483
$ws.Cells.Item(30, 6).Value = 447
$ws.Cells.Item(31, 6).Value = 1184
$ws.Cells.Item(33, 6).Value = 862
$ws.Cells.Item(34, 6).Value = 2578
$ws.Cells.Item(35, 6).Value = 1063
$ws.Cells.Item(36, 6).Value = 576
$ws.Cells.Item(37, 6).Value = 113
$ws.Cells.Item(39, 6).Value = 46
$ws.Cells.Item(41, 6).Value = 817
$ws.Cells.Item(43, 6).Value = 1174

